$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The training session is now split into parts ("Trening"): drop the two
# trailing raw-measurement rows and blank out the per-row Seconds/Velocity/
# Acceleration_SMA measurements that no longer apply once rows are grouped
# by training segment.
$ws.Range("A6:E7").EntireRow.Delete()

# Add the new "Trening" column header, copying the look of the other
# header cells (bold, centered, bordered).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Trening"

# The raw per-row measurements are no longer tracked individually.
$ws.Range("B2:D5").ClearContents()

# Timestamps collapse down to the (shared) session date, stored as a real
# date serial with a date-time number format.
$ws.Range("A2").Value = 45674
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A3").Value = 45674
$ws.Range("A3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A4").Value = 45674
$ws.Range("A4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A5").Value = 45674
$ws.Range("A5").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Velocity_Bin for row 3 moves from 10-15 to 5-10.
$ws.Range("E3").Value = "5-10"

# Fill in which part of the training each row belongs to.
$ws.Range("F2").Value = "Duża Gra"
$ws.Range("F3").Value = "Duża Gra"
$ws.Range("F4").Value = "Mała Gra"
$ws.Range("F5").Value = "Mała Gra"
